$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Class" and "Min_size_mm" columns (old D:E); this shifts
# gamma1/gamma2/Caldens_mn/Caldens_sd left from F:I into D:G and also
# drops the now-unused "Class"/"Min_size_mm" shared strings.
$ws.Range("D1:E1").EntireColumn.Delete()

# Add the two new columns of data (log_G_mn / log_G_sd).
$ws.Range("H1").Value = "log_G_mn"
$ws.Range("I1").Value = "log_G_sd"

$ws.Range("B2").Value = "urchin"
$ws.Range("C2").Value = "urchins, various sp"
$ws.Range("D2").Value = -9.1557912800493906
$ws.Range("E2").Value = 3.1172105405765098
$ws.Range("F2").Value = -0.28768207245178101
$ws.Range("G2").Value = 0.037781937069642203
$ws.Range("H2").Value = 2.0137508156211599
$ws.Range("I2").Value = 0.12278231544543
$ws.Range("B3").Value = "mussel"
$ws.Range("C3").Value = "mussels"
$ws.Range("D3").Value = -6.1239876951112802
$ws.Range("E3").Value = 2.0835640977838201
$ws.Range("F3").Value = -0.28903109543787397
$ws.Range("G3").Value = 0.0625
$ws.Range("H3").Value = 1.3390248208506701
$ws.Range("I3").Value = 0.113096287026453
$ws.Range("B4").Value = "clam"
$ws.Range("C4").Value = "clams, various species"
$ws.Range("D4").Value = -8.6152793977053204
$ws.Range("E4").Value = 2.9069444229797599
$ws.Range("F4").Value = -0.460616762512638
$ws.Range("G4").Value = 0.0625
$ws.Range("H4").Value = 2.5424261328162698
$ws.Range("I4").Value = 0.16728201738732901
$ws.Range("B5").Value = "abalone"
$ws.Range("C5").Value = "abalone, various sp"
$ws.Range("D5").Value = -7.8250485135578201
$ws.Range("E5").Value = 2.68367446900147
$ws.Range("F5").Value = -0.014228674612371201
$ws.Range("G5").Value = 0.0625
$ws.Range("H5").Value = 4.5187521446107199
$ws.Range("I5").Value = 0.179485506757257
$ws.Range("B6").Value = "cancrid_crab"
$ws.Range("C6").Value = "Cancr crabs"
$ws.Range("D6").Value = -9.27390081238131
$ws.Range("E6").Value = 2.96116021712495
$ws.Range("F6").Value = -0.21835233805686899
$ws.Range("G6").Value = 0.0545198284024545
$ws.Range("H6").Value = 3.55296842843984
$ws.Range("I6").Value = 0.19863749715604201
$ws.Range("B7").Value = "kelp_crab"
$ws.Range("C7").Value = "kelp crabs"
$ws.Range("D7").Value = -8.7933118956100103
$ws.Range("E7").Value = 3.05773753996605
$ws.Range("F7").Value = -0.29668997638178901
$ws.Range("G7").Value = 0.0625
$ws.Range("H7").Value = 2.5546121440275198
$ws.Range("I7").Value = 0.167053639049039
$ws.Range("B8").Value = "crab_other"
$ws.Range("C8").Value = "Other crabs"
$ws.Range("D8").Value = -7.2107642201959399
$ws.Range("E8").Value = 2.5270538062199299
$ws.Range("F8").Value = -0.232007386024601
$ws.Range("G8").Value = 0.060578240153726502
$ws.Range("H8").Value = 1.8655089093095001
$ws.Range("I8").Value = 0.15020469041586601
$ws.Range("B9").Value = "snail"
$ws.Range("C9").Value = "snails, various sp"
$ws.Range("D9").Value = -7.9596874196781204
$ws.Range("E9").Value = 2.7738228422176801
$ws.Range("F9").Value = 0.039100135673743902
$ws.Range("G9").Value = 0.0625
$ws.Range("H9").Value = 0.76349440856147699
$ws.Range("I9").Value = 0.14648191676980599
$ws.Range("B10").Value = "star"
$ws.Range("C10").Value = "sea stars"
$ws.Range("D10").Value = -3.0255056410644898
$ws.Range("E10").Value = 1.6308922912286801
$ws.Range("F10").Value = 0.26235164220611601
$ws.Range("G10").Value = 0.051150470442866899
$ws.Range("H10").Value = 3.6245099212635998
$ws.Range("I10").Value = 0.15610395082053299
$ws.Range("B11").Value = "cephalapod"
$ws.Range("C11").Value = "octopus and squid"
$ws.Range("D11").Value = -1.27550564106449
$ws.Range("E11").Value = 1.3127794847699099
$ws.Range("F11").Value = 0.058955356525265998
$ws.Range("G11").Value = 0.061018000103057399
$ws.Range("H11").Value = 4.4069120101129204
$ws.Range("I11").Value = 0.16003677154114301
$ws.Range("B12").Value = "other"
$ws.Range("C12").Value = "chitons, limpets, barnacles,etc."
$ws.Range("D12").Value = -5.2186653100931002
$ws.Range("E12").Value = 1.7983291177557701
$ws.Range("F12").Value = -0.74606822028809305
$ws.Range("G12").Value = 0.0385562741815463
$ws.Range("H12").Value = 1.44661751795251
$ws.Range("I12").Value = 0.100660297301433

# Column width tweaks to match the reformatted table.
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(6).ColumnWidth = 14.666666666666666
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 9.833333333333334

# Update the active selection, as left by the editing session.
$ws.Range("F14").Select()
